# Update the "Estado" (Status) column for the Sobre 3 backlog items
# (rows 19-23) from "Por hacer" to "Hecho".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("B19:B23").Value = "Hecho"

# Move the active selection, matching the saved cursor position.
$ws.Range("D15").Select()
